$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Dadra and Nagar Haveli" (row 9) and "Daman and Diu" (row 10) were merged
# into a single region "Dadra and Nagar Haveli and Daman and Diu". Its
# Population and Density figures are now the sum of the two former regions,
# expressed as formulas; Growth rate / Density Class are kept from the
# original "Dadra and Nagar Haveli" row.
$ws.Cells.Item(9, 1).Value = "Dadra and Nagar Haveli and Daman and Diu"
$ws.Cells.Item(9, 2).Formula = "=343709+243247"
$ws.Cells.Item(9, 3).Formula = "=700+2191"

# Remove the now-redundant "Daman and Diu" row (old row 10), shifting the
# rest of the table up by one row. Only columns A:E are shifted so the
# stray formatted (but empty) F column is left untouched, matching the
# original layout.
$ws.Range("A10:E10").Delete(-4162) # xlShiftUp

# The blank, but formatted, trailing row (F37) survives a plain row
# delete in real Excel even though the used data range shrinks to row 36;
# re-stamp its format so it is retained the same way here.
$ws.Cells.Item(36, 6).Copy() | Out-Null
$ws.Cells.Item(37, 6).PasteSpecial(-4122) | Out-Null # xlPasteFormats
$excel.CutCopyMode = $false

# Re-apply the table's sort (by State Name) now that the table is one row
# shorter.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A36")) | Out-Null
$sortObj.SetRange($ws.Range("A1:E36"))
$sortObj.Header = 1
$sortObj.Apply()

# Column width tweaks: State Name column widened to fit the new longer
# name, and a new width set for column C.
$ws.Columns.Item(1).ColumnWidth = 45.3
$ws.Columns.Item(3).ColumnWidth = 20.9

# Cosmetic: last selected cell.
$ws.Range("A8").Select() | Out-Null
